$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 102449.8
$ws.Range("I40").Value = 127312.375
$ws.Range("J40").Value = 2999.5
$ws.Range("K40").Value = 127312.375
$ws.Range("L40").Value = 2999.5
$ws.Range("M40").Value = -127137.375
$ws.Range("N40").Value = -3349.5
$ws.Range("H76").Value = 5345.6665
$ws.Range("I76").Value = 5026.143
$ws.Range("K76").Value = 5026.143
$ws.Range("M76").Value = -4711.143
$ws.Range("H79").Value = 5345.6665
$ws.Range("I79").Value = 5026.143
$ws.Range("K79").Value = 5026.143
$ws.Range("M79").Value = -3934.143
$ws.Range("H129").Value = 2708.1272
$ws.Range("J129").Value = 940.29785
$ws.Range("L129").Value = 2820.89355
$ws.Range("N129").Value = -12820.89355
$ws.Range("H137").Value = 1763.85
$ws.Range("I137").Value = 1853.8
$ws.Range("K137").Value = 5561.4
$ws.Range("M137").Value = -3011.4
$ws.Range("H138").Value = 3973.8125
$ws.Range("I138").Value = 4012.6
$ws.Range("J138").Value = 3968.2715
$ws.Range("K138").Value = 12037.8
$ws.Range("L138").Value = 11904.8145
$ws.Range("M138").Value = -6897.799999999999
$ws.Range("N138").Value = -22184.8145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 11999
$ws.Range("J9").Value = 11999
$ws.Range("L9").Value = 11999
$ws.Range("N9").Value = -12339
$ws.Range("H20").Value = 11999
$ws.Range("J20").Value = 11999
$ws.Range("L20").Value = 11999
$ws.Range("N20").Value = -12539
$ws.Range("H23").Value = 19876.5
$ws.Range("J23").Value = 9833.333000000001
$ws.Range("L23").Value = 9833.333000000001
$ws.Range("N23").Value = -10351.333
$ws.Range("H32").Value = 50113.188
$ws.Range("I32").Value = 20772.387
$ws.Range("K32").Value = 20772.387
$ws.Range("M32").Value = -20485.387

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1203.6666
$ws.Range("I80").Value = 733.1667
$ws.Range("J80").Value = 1517.3334
$ws.Range("K80").Value = 733.1667
$ws.Range("L80").Value = 1517.3334
$ws.Range("M80").Value = 264.8333
$ws.Range("N80").Value = -3513.3334
$ws.Range("H82").Value = 14209.667
$ws.Range("J82").Value = 28985
$ws.Range("L82").Value = 28985
$ws.Range("N82").Value = -29751
$ws.Range("H83").Value = 1203.6666
$ws.Range("I83").Value = 733.1667
$ws.Range("J83").Value = 1517.3334
$ws.Range("K83").Value = 3665.8335
$ws.Range("L83").Value = 7586.666999999999
$ws.Range("M83").Value = 1326.1665
$ws.Range("N83").Value = -17570.667
$ws.Range("H85").Value = 14209.667
$ws.Range("J85").Value = 28985
$ws.Range("L85").Value = 28985
$ws.Range("N85").Value = -31637
$ws.Range("H86").Value = 71477.125
$ws.Range("I86").Value = 87613.62
$ws.Range("J86").Value = 1552.3334
$ws.Range("K86").Value = 87613.62
$ws.Range("L86").Value = 1552.3334
$ws.Range("M86").Value = -86490.62
$ws.Range("N86").Value = -3798.3334
$ws.Range("H89").Value = 71477.125
$ws.Range("I89").Value = 87613.62
$ws.Range("J89").Value = 1552.3334
$ws.Range("K89").Value = 438068.1
$ws.Range("L89").Value = 7761.666999999999
$ws.Range("M89").Value = -432452.1
$ws.Range("N89").Value = -18993.667
$ws.Range("H107").Value = 142925940
$ws.Range("I107").Value = 250114930
$ws.Range("J107").Value = 7260
$ws.Range("K107").Value = 250114930
$ws.Range("L107").Value = 7260
$ws.Range("M107").Value = -250113010
$ws.Range("N107").Value = -11100
$ws.Range("H134").Value = 3789.7646
$ws.Range("I134").Value = 3647.9124
$ws.Range("J134").Value = 4524.8184
$ws.Range("K134").Value = 10943.7372
$ws.Range("L134").Value = 13574.4552
$ws.Range("M134").Value = -8408.7372
$ws.Range("N134").Value = -18644.4552

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9492
$ws.Range("J50").Value = 9492
$ws.Range("L50").Value = 9492
$ws.Range("N50").Value = -10742
$ws.Range("H51").Value = 7826.125
$ws.Range("J51").Value = 7826.125
$ws.Range("L51").Value = 7826.125
$ws.Range("N51").Value = -9298.125
$ws.Range("H60").Value = 15103
$ws.Range("J60").Value = 15103
$ws.Range("L60").Value = 15103
$ws.Range("N60").Value = -16125
$ws.Range("H61").Value = 7826.125
$ws.Range("J61").Value = 7826.125
$ws.Range("L61").Value = 7826.125
$ws.Range("N61").Value = -8522.125
$ws.Range("H68").Value = 19729.54
$ws.Range("J68").Value = 19729.54
$ws.Range("L68").Value = 19729.54
$ws.Range("N68").Value = -21227.54
$ws.Range("H71").Value = 19729.54
$ws.Range("J71").Value = 19729.54
$ws.Range("L71").Value = 59188.62
$ws.Range("N71").Value = -66676.62
$ws.Range("H74").Value = 38986.668
$ws.Range("J74").Value = 38986.668
$ws.Range("L74").Value = 38986.668
$ws.Range("N74").Value = -40734.668
$ws.Range("H77").Value = 38986.668
$ws.Range("J77").Value = 38986.668
$ws.Range("L77").Value = 116960.004
$ws.Range("N77").Value = -125696.004
$ws.Range("H99").Value = 15111.75
$ws.Range("I99").Value = 3690
$ws.Range("J99").Value = 18919
$ws.Range("K99").Value = 3690
$ws.Range("L99").Value = 18919
$ws.Range("M99").Value = -2192
$ws.Range("N99").Value = -21915
$ws.Range("H126").Value = 15111.75
$ws.Range("I126").Value = 3690
$ws.Range("J126").Value = 18919
$ws.Range("K126").Value = 11070
$ws.Range("L126").Value = 56757
$ws.Range("M126").Value = -8600
$ws.Range("N126").Value = -61697
$ws.Range("H132").Value = 2425.7334
$ws.Range("I132").Value = 2378.36
$ws.Range("J132").Value = 2662.6
$ws.Range("K132").Value = 7135.08
$ws.Range("L132").Value = 7987.799999999999
$ws.Range("M132").Value = -4605.08
$ws.Range("N132").Value = -13047.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1843.85
$ws.Range("J5").Value = 1666.6666
$ws.Range("L5").Value = 4999.9998
$ws.Range("N5").Value = -5223.9998
$ws.Range("H6").Value = 272
$ws.Range("I6").Value = 44
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 132
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -19
$ws.Range("N6").Value = -1726
$ws.Range("H131").Value = 768989.75
$ws.Range("J131").Value = 805571.8
$ws.Range("L131").Value = 2416715.4
$ws.Range("N131").Value = -2426795.4
$ws.Range("H135").Value = 1843.85
$ws.Range("J135").Value = 1666.6666
$ws.Range("L135").Value = 14999.9994
$ws.Range("N135").Value = -20069.9994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 129031.75
$ws.Range("I70").Value = 203710.8
$ws.Range("K70").Value = 203710.8
$ws.Range("M70").Value = -203440.8
$ws.Range("H73").Value = 129031.75
$ws.Range("I73").Value = 203710.8
$ws.Range("K73").Value = 203710.8
$ws.Range("M73").Value = -202774.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8546.666999999999
$ws.Range("I46").Value = 9340
$ws.Range("J46").Value = 8150
$ws.Range("K46").Value = 9340
$ws.Range("L46").Value = 8150
$ws.Range("M46").Value = -9152
$ws.Range("N46").Value = -8526
$ws.Range("H100").Value = 3231.25
$ws.Range("I100").Value = 2422.5
$ws.Range("K100").Value = 2422.5
$ws.Range("M100").Value = -1881.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3278.7795
$ws.Range("I136").Value = 4632.923
$ws.Range("J136").Value = 2211.879
$ws.Range("K136").Value = 13898.769
$ws.Range("L136").Value = 6635.637
$ws.Range("M136").Value = -11348.769
$ws.Range("N136").Value = -11735.637
